# Utah roster: Lauri Markkanen and Jordan Clarkson swap places in the
# table (row 3 <-> row 4). Column A ("No.", the row index) stays put;
# everything else (jersey #, name, position, height, weight, birth date,
# country, experience, college, bbref url) moves with the player.
#
# Use Copy/PasteSpecial (rather than reading/writing .Value2 directly)
# so that each cell's original type/formatting - e.g. the "Exp" column
# values that are stored as text ("5", "8") rather than numbers - is
# preserved exactly as it swaps rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash row 3 in an unused scratch row far below the table.
$ws.Range("B3:K3").Copy()
$ws.Range("B100:K100").PasteSpecial()

# Row 4 -> row 3
$ws.Range("B4:K4").Copy()
$ws.Range("B3:K3").PasteSpecial()

# Scratch (old row 3) -> row 4
$ws.Range("B100:K100").Copy()
$ws.Range("B4:K4").PasteSpecial()

# Clean up the scratch row and the marching-ants clipboard marquee.
$ws.Range("B100:K100").Clear()
$excel.CutCopyMode = $false
